# Update "想去人数" (want-to-go count) figures in column F
# for the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 8403
    "F3"  = 7969
    "F10" = 180
    "F11" = 236
    "F13" = 140
    "F14" = 2024
    "F19" = 134
    "F20" = 33
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
